$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.190.58'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.569.11'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.68'
$ws.Range("E5").Value = '  +1.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.493'
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  -0.64%  '
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0600'
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").Value = '1.792.34'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '1.546.04'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.79'
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '27.173.86'
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.30'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.42'
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '215.86'
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.24'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.90'
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.66'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.09'
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.107'
$ws.Range("E28").Value = '  +1.87%  '
$ws.Range("E30").Value = '  +2.46%  '
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.25'
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("D34").Value = '1.448.14'
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("E35").Value = '  +5.04%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.535'
$ws.Range("E39").Value = '  +0.64%  '
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.73'
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").Value = '1.706.24'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("E49").Value = '  +3.18%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("E51").Value = '  +0.10%  '
